# Update symbol list (price / volume columns) to reflect the latest
# crypto pricing snapshot, as produced by the scheduled GitHub Actions
# refresh job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing Excel to treat it as
# plain text (so values like "0.03000" or "0.09%" keep their exact
# textual formatting instead of being silently reinterpreted as
# numbers/percentages). We do this by assigning a formula that evaluates
# to the literal text, then flattening that formula down to its static
# result via copy / paste-values, which leaves the cell's style/format
# untouched.
function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$TextValue
    )

    $cell = $ws.Range($CellRef)
    $escaped = $TextValue.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

Set-TextValue "D2" "298.90"
Set-TextValue "E2" "-0.63%"

Set-TextValue "D3" "31.58"
Set-TextValue "E3" "0.72%"

Set-TextValue "D4" "5.149"
Set-TextValue "E4" "0.23%"

Set-TextValue "D5" "0.08095"
Set-TextValue "E5" "9.50%"

Set-TextValue "D6" "2.496"
Set-TextValue "E6" "14.75%"

Set-TextValue "D7" "7.800"
Set-TextValue "E7" "-1.57%"

Set-TextValue "D8" "3.911"
Set-TextValue "E8" "2.28%"

Set-TextValue "D9" "0.9296"
Set-TextValue "E9" "1.20%"

Set-TextValue "D10" "0.1761"
Set-TextValue "E10" "3.05%"

Set-TextValue "D11" "0.07416"
Set-TextValue "E11" "-2.08%"

Set-TextValue "D12" "0.08882"
Set-TextValue "E12" "8.95%"

Set-TextValue "D13" "0.02998"
Set-TextValue "E13" "-0.49%"

Set-TextValue "D14" "0.09995"
Set-TextValue "E14" "0.67%"

Set-TextValue "E15" "1.17%"

Set-TextValue "D16" "0.005957"
Set-TextValue "E16" "-1.92%"

Set-TextValue "D17" "3.530"
Set-TextValue "E17" "1.64%"

Set-TextValue "D18" "2.290"
Set-TextValue "E18" "2.96%"

Set-TextValue "E19" "0.35%"

Set-TextValue "E20" "1.49%"

Set-TextValue "D21" "4.163"
Set-TextValue "E21" "-10.61%"

Set-TextValue "D22" "0.1679"
Set-TextValue "E22" "7.23%"

Set-TextValue "E23" "-0.36%"

Set-TextValue "D24" "0.001241"
Set-TextValue "E24" "1.28%"

Set-TextValue "D25" "0.004535"
Set-TextValue "E25" "1.35%"

Set-TextValue "D26" "0.0001200"
Set-TextValue "E26" "-7.68%"

Set-TextValue "D27" "0.0003407"
Set-TextValue "E27" "-0.53%"

Set-TextValue "D39" "0.01754"
Set-TextValue "E39" "0.60%"

Set-TextValue "D40" "0.04591"
Set-TextValue "E40" "1.69%"

Set-TextValue "D41" "0.006919"
Set-TextValue "E41" "-5.36%"

Set-TextValue "E42" "1.86%"

Set-TextValue "E43" "-1.77%"

Set-TextValue "D44" "0.01029"
Set-TextValue "E44" "-3.23%"

Set-TextValue "D45" "0.00006145"
Set-TextValue "E45" "-1.99%"

Set-TextValue "D46" "0.00000000749"
Set-TextValue "E46" "-0.11%"

Set-TextValue "D47" "0.008392"
Set-TextValue "E47" "-16.01%"

Set-TextValue "D48" "0.7484"
Set-TextValue "E48" "-7.43%"

Set-TextValue "D49" "0.00002098"
Set-TextValue "E49" "-0.11%"

Set-TextValue "D50" "0.0001998"
Set-TextValue "E50" "-0.04%"

$excel.CutCopyMode = 0
